$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 151, pushing the existing rows 151-264 down to 152-265.
$ws.Rows(151).Insert()

# Populate the newly inserted row 151 with the new record's data.
$ws.Range("A151").Value = 3
$ws.Range("B151").Value = "Femacal de La Calera"
$ws.Range("C151").Value = "Coquimbo"
$ws.Range("D151").Value = 44651
$ws.Range("E151").Value = 5
$ws.Range("F151").Value = 100112001
$ws.Range("G151").Value = "Berenjena"
$ws.Range("H151").Value = "Sin especificar"
$ws.Range("I151").Value = "Primera"
$ws.Range("J151").Value = 115
$ws.Range("K151").Value = 9000
$ws.Range("L151").Value = 9500
$ws.Range("M151").Value = 9261
$ws.Range("N151").Value = '$/caja 60 unidades'
$ws.Range("O151").Value = "Región de Arica y Parinacota"
$ws.Range("P151").Value = 154
$ws.Range("Q151").Value = 60
$ws.Range("R151").Value = "Hortaliza"
